$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1168.3494
$ws.Range("J17").Value = 1160.2839
$ws.Range("L17").Value = 3480.8517
$ws.Range("N17").Value = -3816.8517
$ws.Range("H51").Value = 3972.8667
$ws.Range("J51").Value = 4012.375
$ws.Range("L51").Value = 4012.375
$ws.Range("N51").Value = -4980.375
$ws.Range("H100").Value = 1744.3448
$ws.Range("I100").Value = 1497.5
$ws.Range("J100").Value = 1918.5883
$ws.Range("K100").Value = 1497.5
$ws.Range("L100").Value = 1918.5883
$ws.Range("M100").Value = -956.5
$ws.Range("N100").Value = -3000.5883
$ws.Range("H135").Value = 873.0172
$ws.Range("I135").Value = 558.2593000000001
$ws.Range("K135").Value = 5024.3337
$ws.Range("M135").Value = -2489.3337
$ws.Range("H137").Value = 4491.96
$ws.Range("J137").Value = 2364.2593
$ws.Range("L137").Value = 7092.777900000001
$ws.Range("N137").Value = -12192.7779
$ws.Range("H138").Value = 16164726
$ws.Range("I138").Value = 32259802
$ws.Range("K138").Value = 96779406
$ws.Range("M138").Value = -96774266
$ws.Range("H141").Value = 6485.6045
$ws.Range("I141").Value = 3496.8684
$ws.Range("K141").Value = 10490.6052
$ws.Range("M141").Value = -5310.6052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4678.226
$ws.Range("I32").Value = 4137.5874
$ws.Range("K32").Value = 4137.5874
$ws.Range("M32").Value = -3850.5874
$ws.Range("H74").Value = 1459.4783
$ws.Range("I74").Value = 1439.7368
$ws.Range("J74").Value = 1553.25
$ws.Range("K74").Value = 1439.7368
$ws.Range("L74").Value = 1553.25
$ws.Range("M74").Value = -565.7367999999999
$ws.Range("N74").Value = -3301.25
$ws.Range("H77").Value = 1459.4783
$ws.Range("I77").Value = 1439.7368
$ws.Range("J77").Value = 1553.25
$ws.Range("K77").Value = 7198.683999999999
$ws.Range("L77").Value = 7766.25
$ws.Range("M77").Value = -2830.683999999999
$ws.Range("N77").Value = -16502.25
$ws.Range("H97").Value = 38430.37
$ws.Range("I97").Value = 1446.9231
$ws.Range("K97").Value = 1446.9231
$ws.Range("M97").Value = -950.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1042
$ws.Range("I94").Value = 996.8333
$ws.Range("J94").Value = 1109.75
$ws.Range("K94").Value = 996.8333
$ws.Range("L94").Value = 1109.75
$ws.Range("M94").Value = -545.8333
$ws.Range("N94").Value = -2011.75
$ws.Range("H99").Value = 1405.9474
$ws.Range("I99").Value = 1261.8889
$ws.Range("K99").Value = 1261.8889
$ws.Range("M99").Value = 236.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3090.7878
$ws.Range("I16").Value = 2483.5454
$ws.Range("J16").Value = 4305.273
$ws.Range("K16").Value = 2483.5454
$ws.Range("L16").Value = 4305.273
$ws.Range("M16").Value = -2196.5454
$ws.Range("N16").Value = -4879.273
$ws.Range("H31").Value = 2031.7046
$ws.Range("I31").Value = 1791.5938
$ws.Range("K31").Value = 1791.5938
$ws.Range("M31").Value = -1496.5938
$ws.Range("H34").Value = 2031.7046
$ws.Range("I34").Value = 1791.5938
$ws.Range("K34").Value = 1791.5938
$ws.Range("M34").Value = -1589.5938
$ws.Range("H107").Value = 1060.9166
$ws.Range("I107").Value = 1106.5714
$ws.Range("J107").Value = 741.3333
$ws.Range("K107").Value = 1106.5714
$ws.Range("L107").Value = 741.3333
$ws.Range("M107").Value = 813.4286
$ws.Range("N107").Value = -4581.3333
$ws.Range("H113").Value = 3090.7878
$ws.Range("I113").Value = 2483.5454
$ws.Range("J113").Value = 4305.273
$ws.Range("K113").Value = 2483.5454
$ws.Range("L113").Value = 4305.273
$ws.Range("M113").Value = -313.5454
$ws.Range("N113").Value = -8645.273000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1186.4324
$ws.Range("I5").Value = 1026.3103
$ws.Range("K5").Value = 3078.9309
$ws.Range("M5").Value = -2966.9309
$ws.Range("H92").Value = 667.1111
$ws.Range("I92").Value = 721
$ws.Range("J92").Value = 646.38464
$ws.Range("K92").Value = 2163
$ws.Range("L92").Value = 1939.15392
$ws.Range("M92").Value = -915
$ws.Range("N92").Value = -4435.15392
$ws.Range("H121").Value = 1914.4445
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 2103.75
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 6311.25
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -8931.25
$ws.Range("H134").Value = 2489.261
$ws.Range("I134").Value = 1960.6842
$ws.Range("K134").Value = 5882.0526
$ws.Range("M134").Value = -812.0526
$ws.Range("H135").Value = 1186.4324
$ws.Range("I135").Value = 1026.3103
$ws.Range("K135").Value = 9236.792700000002
$ws.Range("M135").Value = -6701.792700000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13533036
$ws.Range("J80").Value = 4460.4443
$ws.Range("L80").Value = 4460.4443
$ws.Range("N80").Value = -6456.4443
$ws.Range("H83").Value = 13533036
$ws.Range("J83").Value = 4460.4443
$ws.Range("L83").Value = 22302.2215
$ws.Range("N83").Value = -32286.2215
$ws.Range("H132").Value = 5298.7017
$ws.Range("I132").Value = 3887.9534
$ws.Range("K132").Value = 11663.8602
$ws.Range("M132").Value = -9133.860199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2303.5
$ws.Range("I7").Value = 2319.1428
$ws.Range("K7").Value = 2319.1428
$ws.Range("M7").Value = -2207.1428
$ws.Range("H93").Value = 2656.9375
$ws.Range("I93").Value = 2608.6
$ws.Range("K93").Value = 2608.6
$ws.Range("M93").Value = -1360.6
$ws.Range("H126").Value = 2303.5
$ws.Range("I126").Value = 2319.1428
$ws.Range("K126").Value = 6957.428400000001
$ws.Range("M126").Value = -4487.428400000001
$ws.Range("H132").Value = 25401.334
$ws.Range("I132").Value = 30221.955
$ws.Range("K132").Value = 90665.86500000001
$ws.Range("M132").Value = -88135.86500000001
$ws.Range("H136").Value = 7504173
$ws.Range("I136").Value = 15001652
$ws.Range("J136").Value = 6694.5835
$ws.Range("K136").Value = 45004956
$ws.Range("L136").Value = 20083.7505
$ws.Range("M136").Value = -45002406
$ws.Range("N136").Value = -25183.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9061.9375
$ws.Range("I62").Value = 5553.2104
$ws.Range("J62").Value = 14190.077
$ws.Range("K62").Value = 5553.2104
$ws.Range("L62").Value = 14190.077
$ws.Range("M62").Value = -4929.2104
$ws.Range("N62").Value = -15438.077
$ws.Range("H65").Value = 9061.9375
$ws.Range("I65").Value = 5553.2104
$ws.Range("J65").Value = 14190.077
$ws.Range("K65").Value = 27766.052
$ws.Range("L65").Value = 70950.38499999999
$ws.Range("M65").Value = -24646.052
$ws.Range("N65").Value = -77190.38499999999
$ws.Range("H122").Value = 4797.95
$ws.Range("I122").Value = 3154.25
$ws.Range("K122").Value = 9462.75
$ws.Range("M122").Value = -7012.75
$ws.Range("H126").Value = 9618098
$ws.Range("I126").Value = 9618098
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 28854294
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -28851824
$ws.Range("N126").ClearContents() | Out-Null
$ws.Range("H132").Value = 3225.1892
$ws.Range("I132").Value = 2494.9395
$ws.Range("J132").Value = 9249.75
$ws.Range("K132").Value = 7484.818499999999
$ws.Range("L132").Value = 27749.25
$ws.Range("M132").Value = -4954.818499999999
$ws.Range("N132").Value = -32809.25
$ws.Range("H136").Value = 1571.2639
$ws.Range("I136").Value = 1682.0754
$ws.Range("K136").Value = 5046.2262
$ws.Range("M136").Value = -2496.2262
